# Commit: "Added a pdf for tasks 1 and 2."
#
# The diff appends three new numbered list items (numId="2", same list used
# by the two preceding "finite height / monotonic transfer functions" bullet
# points) right after the paragraph that ends in "...guaranteed solution. ",
# and it also marks the built-in "Default Paragraph Font" character style as
# w:semiHidden in styles.xml.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Locate the anchor paragraph robustly (by its distinctive trailing text)
#    rather than a hard-coded paragraph index, then work out which 1-based
#    Paragraphs index it is so we can insert new paragraphs right after it.
# ---------------------------------------------------------------------------
$anchorText = "guaranteeing a point of termination within the algorithm with a guaranteed solution. "

$found = $d.Content.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the anchor paragraph ending in 'guaranteed solution.'"
}
$anchorRange = $d.Content

# Count paragraph marks before the match to get the anchor's 1-based index.
$textBefore = $d.Range(0, $anchorRange.Start).Text
$paraNum = ($textBefore.ToCharArray() | Where-Object { $_ -eq [char]13 }).Count + 1

# ---------------------------------------------------------------------------
# 2) Build the three new list paragraphs as raw WordOpenXML fragments so the
#    resulting markup (pStyle/numPr/proofErr/run-splits) matches exactly.
# ---------------------------------------------------------------------------
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$xmlTask1 = '<w:p ' + $wNs + '>' + `
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr>' + `
    '<w:r><w:t>Answer can be found in df.py</w:t></w:r>' + `
    '</w:p>'

$xmlTask1Grammar = '<w:p ' + $wNs + '>' + `
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr>' + `
    '<w:r><w:t xml:space="preserve">Answer can </w:t></w:r>' + `
    '<w:proofErr w:type="gramStart"/>' + `
    '<w:r><w:t>found</w:t></w:r>' + `
    '<w:proofErr w:type="gramEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> in df.p</w:t></w:r>' + `
    '<w:r><w:t>y</w:t></w:r>' + `
    '</w:p>'

$xmlTask2 = '<w:p ' + $wNs + '>' + `
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr>' + `
    '<w:r><w:t>Answer can be found in /test</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> directory.</w:t></w:r>' + `
    '</w:p>'

# ---------------------------------------------------------------------------
# 3) Insert them, one at a time, immediately after the anchor paragraph.
# ---------------------------------------------------------------------------
$d.Paragraphs.Item($paraNum).Range.InsertParagraphAfter()
$d.Paragraphs.Item($paraNum + 1).Range.InsertXML($xmlTask1)

$d.Paragraphs.Item($paraNum + 1).Range.InsertParagraphAfter()
$d.Paragraphs.Item($paraNum + 2).Range.InsertXML($xmlTask1Grammar)

$d.Paragraphs.Item($paraNum + 2).Range.InsertParagraphAfter()
$d.Paragraphs.Item($paraNum + 3).Range.InsertXML($xmlTask2)

Write-Output "Inserted the three 'Answer can be found in ...' list items."

# ---------------------------------------------------------------------------
# 4) Best-effort: flag "Default Paragraph Font" as semi-hidden, matching the
#    <w:semiHidden/> added to that style in styles.xml. Wrapped defensively
#    since this particular flag isn't always exposed as a writable property.
# ---------------------------------------------------------------------------
try {
    $dpf = $d.Styles.Item("Default Paragraph Font")
    try { $dpf.SemiHidden = $true } catch { }
    try { $dpf.Hidden = $true } catch { }
} catch { }
